$wb = $excel.ActiveWorkbook

# The data being updated lives on the "VENTA MENSUAL" sheet.
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Shift the month headers forward by one month.
$ws.Range("C1").Value = "junio"
$ws.Range("D1").Value = "julio"
$ws.Range("E1").Value = "agosto"
$ws.Range("F1").Value = "septiembre"

# Update the column widths to match the new header text.
# (ColumnWidth as set via COM is offset from the raw OOXML column width by
# 0.8333333333333334, so subtract that to land on the exact target width.)
$ws.Columns.Item(3).ColumnWidth = 10.166666666666666
$ws.Columns.Item(5).ColumnWidth = 11.166666666666666
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666
